$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 9 new rows (rows 22-30), continuing the existing table pattern.
$regCenterIds = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$machineIds   = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)

for ($i = 0; $i -lt $regCenterIds.Count; $i++) {
    $row = 22 + $i

    $ws.Cells.Item($row, 1).Value = $regCenterIds[$i]   # A - regcntr_id
    $ws.Cells.Item($row, 2).Value = $machineIds[$i]     # B - machine_id
    $ws.Cells.Item($row, 3).Value = "eng"                # C - lang_code
    $ws.Cells.Item($row, 4).Value = $true                # D - is_active
    $ws.Cells.Item($row, 5).Value = "superadmin"          # E - cr_by
    $ws.Cells.Item($row, 6).Value = "now()"               # F - cr_dtimes
}

# Update the sheet view: select the rows below the data (mirrors the
# "select to end" state left behind after entering the new rows in Excel).
[void]$ws.Range("A31:XFD1048576").Select()

# Add a basic page setup (portrait orientation) to the sheet, as recorded
# in the workbook after the edit.
$ws.PageSetup.Orientation = 1
